$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post that used to occupy row 780 ("「忍耐の後に、待望の美しいこと」")
# was removed from the source data. Deleting the entire row shifts every
# subsequent row (781-839) up by one, which matches the new dimension
# A1:C838 and the renumbered rows shown in the diff.
$ws.Rows("780:780").Delete()
